# add endpoint for User Question
# Appends a new question/answer row ("ما اسم مدير يونا" / "عاطف") to Sheet1,
# mirroring the existing question/answer rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "ما اسم مدير يونا"
$ws.Range("B5").Value = "عاطف"

# Match the author's last selection recorded in the sheet (B5).
$ws.Range("B5").Select()
